$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two "assembly id" entries (shown in column G on the section header rows)
$ws.Range("G2").Value = "FR_A0500"
$ws.Range("G5").Value = "FR_A0700"

# Update the active selection to match the saved view state
$ws.Range("G6").Select()
